$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column-D price cells to remain Text so numeric-looking values
# (e.g. "1.000", "0.6930", "29.441.39") keep their exact literal formatting
# instead of being auto-coerced into floating point numbers by the Value setter.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '29.441.39'
$ws.Cells.Item(2, 5).Value = '  +1.81%  '
$ws.Cells.Item(3, 4).Value = '1.855.51'
$ws.Cells.Item(3, 5).Value = '  +1.21%  '
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).Value = '244.66'
$ws.Cells.Item(5, 5).Value = '  -0.13%  '
$ws.Cells.Item(6, 4).Value = '0.6930'
$ws.Cells.Item(6, 5).Value = '  +0.32%  '
$ws.Cells.Item(7, 4).Value = '1.001'
$ws.Cells.Item(7, 5).Value = '  +0.11%  '
$ws.Cells.Item(8, 4).Value = '0.07653'
$ws.Cells.Item(8, 5).Value = '  -0.38%  '
$ws.Cells.Item(9, 4).Value = '0.3054'
$ws.Cells.Item(9, 5).Value = '  +0.05%  '
$ws.Cells.Item(10, 4).Value = '23.52'
$ws.Cells.Item(10, 5).Value = '  -0.02%  '
$ws.Cells.Item(11, 5).Value = '  -0.68%  '
$ws.Cells.Item(12, 4).Value = '1.851.68'
$ws.Cells.Item(12, 5).Value = '  +1.14%  '
$ws.Cells.Item(13, 4).Value = '5.126'
$ws.Cells.Item(13, 5).Value = '  +1.07%  '
$ws.Cells.Item(14, 4).Value = '0.6894'
$ws.Cells.Item(14, 5).Value = '  +1.36%  '
$ws.Cells.Item(15, 4).Value = '90.38'
$ws.Cells.Item(15, 5).Value = '  -0.17%  '
$ws.Cells.Item(16, 4).Value = '6.424'
$ws.Cells.Item(16, 5).Value = '  -0.17%  '
$ws.Cells.Item(17, 4).Value = '29.355.37'
$ws.Cells.Item(17, 5).Value = '  +1.49%  '
$ws.Cells.Item(18, 4).Value = '0.000008249'
$ws.Cells.Item(18, 5).Value = '  -1.18%  '
$ws.Cells.Item(19, 4).Value = '2.107.98'
$ws.Cells.Item(19, 5).Value = '  +1.29%  '
$ws.Cells.Item(20, 4).Value = '237.27'
$ws.Cells.Item(20, 5).Value = '  -2.44%  '
$ws.Cells.Item(21, 4).Value = '12.70'
$ws.Cells.Item(21, 5).Value = '  +0.12%  '
$ws.Cells.Item(22, 4).Value = '1.000'
$ws.Cells.Item(22, 5).Value = '  +0.07%  '
$ws.Cells.Item(23, 4).Value = '7.620'
$ws.Cells.Item(23, 5).Value = '  +1.89%  '
$ws.Cells.Item(24, 4).Value = '1.001'
$ws.Cells.Item(24, 5).Value = '  +0.16%  '
$ws.Cells.Item(25, 4).Value = '0.1489'
$ws.Cells.Item(25, 5).Value = '  +1.49%  '
$ws.Cells.Item(26, 4).Value = '8.889'
$ws.Cells.Item(26, 5).Value = '  +1.04%  '
$ws.Cells.Item(27, 4).Value = '159.62'
$ws.Cells.Item(27, 5).Value = '  -1.27%  '
$ws.Cells.Item(28, 4).Value = '18.23'
$ws.Cells.Item(28, 5).Value = '  +0.20%  '
$ws.Cells.Item(29, 4).Value = '1.533'
$ws.Cells.Item(29, 5).Value = '  -1.39%  '
$ws.Cells.Item(30, 4).Value = '4.249'
$ws.Cells.Item(30, 5).Value = '  +0.85%  '
$ws.Cells.Item(31, 4).Value = '4.147'
$ws.Cells.Item(31, 5).Value = '  -0.19%  '
$ws.Cells.Item(32, 4).Value = '1.192'
$ws.Cells.Item(32, 5).Value = '  +1.31%  '
$ws.Cells.Item(33, 4).Value = '0.05106'
$ws.Cells.Item(33, 5).Value = '  -0.33%  '
$ws.Cells.Item(34, 4).Value = '0.7669'
$ws.Cells.Item(34, 5).Value = '  +0.22%  '
$ws.Cells.Item(35, 4).Value = '1.883'
$ws.Cells.Item(35, 5).Value = '  +1.97%  '
$ws.Cells.Item(36, 5).Value = '  +0.19%  '
$ws.Cells.Item(37, 5).Value = '  +0.25%  '
$ws.Cells.Item(38, 4).Value = '1.327.55'
$ws.Cells.Item(38, 5).Value = '  +7.79%  '
$ws.Cells.Item(39, 4).Value = '0.01858'
$ws.Cells.Item(39, 5).Value = '  +0.85%  '
$ws.Cells.Item(40, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(40, 4).Value = '0.9797'
$ws.Cells.Item(40, 5).Value = '  +6.42%  '
$ws.Cells.Item(41, 2).Value = 'MXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(41, 4).Value = '2.721'
$ws.Cells.Item(41, 5).Value = '  +1.08%  '
$ws.Cells.Item(42, 4).Value = '105.98'
$ws.Cells.Item(42, 5).Value = '  -2.26%  '
$ws.Cells.Item(43, 4).Value = '5.828'
$ws.Cells.Item(43, 5).Value = '  -0.25%  '
$ws.Cells.Item(44, 4).Value = '1.001'
$ws.Cells.Item(44, 5).Value = '  +0.11%  '
$ws.Cells.Item(45, 4).Value = '0.00000000127'
$ws.Cells.Item(45, 5).Value = '  +4.26%  '
$ws.Cells.Item(46, 4).Value = '9.773'
$ws.Cells.Item(46, 5).Value = '  +2.24%  '
$ws.Cells.Item(47, 4).Value = '2.003.54'
$ws.Cells.Item(47, 5).Value = '  +1.02%  '
$ws.Cells.Item(48, 4).Value = '0.5220'
$ws.Cells.Item(48, 5).Value = '  +0.99%  '
$ws.Cells.Item(49, 5).Value = '  +1.48%  '
$ws.Cells.Item(50, 4).Value = '62.82'
$ws.Cells.Item(50, 5).Value = '  -1.65%  '
$ws.Cells.Item(51, 4).Value = '6.954'
$ws.Cells.Item(51, 5).Value = '  +0.38%  '
